# ---------------------------------------------------------------------------
# Add "_dev_info" column to the "fwk_content" sheet of the ANSSI MonAideCyber
# questionnaire workbook.
#
# This mirrors the upstream commit:
#   "[Excel + Script] Add "_dev_info" column"
#   Informational column in order to ease future updates of the questionnaire
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# "fwk_content" is the 3rd sheet (library_meta, fwk_meta, fwk_content, ...)
$ws = $wb.Worksheets.Item("fwk_content")

# --- New column L: "_dev_info" header + one "Question N" label per question row
$ws.Range("L1").NumberFormat = "@"
$ws.Range("L1").Value = "_dev_info"
$ws.Range("L2").NumberFormat = "@"
$ws.Range("L3").NumberFormat = "@"
$ws.Range("L3").Value = "Question 1"
$ws.Range("L4").NumberFormat = "@"
$ws.Range("L4").Value = "Question 2"
$ws.Range("L5").NumberFormat = "@"
$ws.Range("L5").Value = "Question 3"
$ws.Range("L6").NumberFormat = "@"
$ws.Range("L6").Value = "Question 4"
$ws.Range("L7").NumberFormat = "@"
$ws.Range("L7").Value = "Question 5"
$ws.Range("L8").NumberFormat = "@"
$ws.Range("L8").Value = "Question 6"
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "Question 7"
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "Question 8"
$ws.Range("L11").NumberFormat = "@"
$ws.Range("L11").Value = "Question 8.1"
$ws.Range("L12").NumberFormat = "@"
$ws.Range("L12").Value = "Question 8.2"
$ws.Range("L13").NumberFormat = "@"
$ws.Range("L14").NumberFormat = "@"
$ws.Range("L14").Value = "Question 9"
$ws.Range("L15").NumberFormat = "@"
$ws.Range("L15").Value = "Question 10"
$ws.Range("L16").NumberFormat = "@"
$ws.Range("L16").Value = "Question 11"
$ws.Range("L17").NumberFormat = "@"
$ws.Range("L17").Value = "Question 12"
$ws.Range("L18").NumberFormat = "@"
$ws.Range("L18").Value = "Question 13"
$ws.Range("L19").NumberFormat = "@"
$ws.Range("L19").Value = "Question 14"
$ws.Range("L20").NumberFormat = "@"
$ws.Range("L21").NumberFormat = "@"
$ws.Range("L21").Value = "Question 15"
$ws.Range("L22").NumberFormat = "@"
$ws.Range("L22").Value = "Question 16"
$ws.Range("L23").NumberFormat = "@"
$ws.Range("L23").Value = "Question 17"
$ws.Range("L24").NumberFormat = "@"
$ws.Range("L24").Value = "Question 18"
$ws.Range("L25").NumberFormat = "@"
$ws.Range("L25").Value = "Question 19"
$ws.Range("L26").NumberFormat = "@"
$ws.Range("L26").Value = "Question 20"
$ws.Range("L27").NumberFormat = "@"
$ws.Range("L27").Value = "Question 21"
$ws.Range("L28").NumberFormat = "@"
$ws.Range("L28").Value = "Question 22"
$ws.Range("L29").NumberFormat = "@"
$ws.Range("L29").Value = "Question 23"
$ws.Range("L30").NumberFormat = "@"
$ws.Range("L30").Value = "Question 24"
$ws.Range("L31").NumberFormat = "@"
$ws.Range("L31").Value = "Question 25"
$ws.Range("L32").NumberFormat = "@"
$ws.Range("L33").NumberFormat = "@"
$ws.Range("L33").Value = "Question 26"
$ws.Range("L34").NumberFormat = "@"
$ws.Range("L34").Value = "Question 27"
$ws.Range("L35").NumberFormat = "@"
$ws.Range("L35").Value = "Question 28"
$ws.Range("L36").NumberFormat = "@"
$ws.Range("L36").Value = "Question 29"
$ws.Range("L37").NumberFormat = "@"
$ws.Range("L37").Value = "Question 30"
$ws.Range("L38").NumberFormat = "@"
$ws.Range("L38").Value = "Question 31"
$ws.Range("L39").NumberFormat = "@"
$ws.Range("L39").Value = "Question 32"
$ws.Range("L40").NumberFormat = "@"
$ws.Range("L41").NumberFormat = "@"
$ws.Range("L41").Value = "Question 33"
$ws.Range("L42").NumberFormat = "@"
$ws.Range("L42").Value = "Question 33.1"
$ws.Range("L43").NumberFormat = "@"
$ws.Range("L43").Value = "Question 33.2"
$ws.Range("L44").NumberFormat = "@"
$ws.Range("L44").Value = "Question 34"
$ws.Range("L45").NumberFormat = "@"
$ws.Range("L45").Value = "Question 35"
$ws.Range("L46").NumberFormat = "@"
$ws.Range("L46").Value = "Question 36"
$ws.Range("L47").NumberFormat = "@"
$ws.Range("L47").Value = "Question 37"
$ws.Range("L48").NumberFormat = "@"
$ws.Range("L48").Value = "Question 38"
$ws.Range("L49").NumberFormat = "@"
$ws.Range("L49").Value = "Question 39"
$ws.Range("L50").NumberFormat = "@"
$ws.Range("L51").NumberFormat = "@"
$ws.Range("L51").Value = "Question 40"
$ws.Range("L52").NumberFormat = "@"
$ws.Range("L52").Value = "Question 41"
$ws.Range("L53").NumberFormat = "@"
$ws.Range("L53").Value = "Question 42"
$ws.Range("L54").NumberFormat = "@"
$ws.Range("L55").NumberFormat = "@"
$ws.Range("L55").Value = "Question 43"
$ws.Range("L56").NumberFormat = "@"
$ws.Range("L56").Value = "Question 44"
$ws.Range("L57").NumberFormat = "@"
$ws.Range("L57").Value = "Question 44.1"
$ws.Range("L58").NumberFormat = "@"
$ws.Range("L58").Value = "Question 44.2"
$ws.Range("L59").NumberFormat = "@"
$ws.Range("L59").Value = "Question 44.3"
$ws.Range("L60").NumberFormat = "@"
$ws.Range("L60").Value = "Question 44.4"
$ws.Range("L61").NumberFormat = "@"
$ws.Range("L61").Value = "Question 45"

# Give the new column roughly the same width Excel would pick when auto-fitting
# the "Question N" labels.
$ws.Columns.Item(12).ColumnWidth = 13.2

# The author's last save left "fwk_content" as the active/selected tab
# (instead of "library_meta").
$ws.Activate() | Out-Null
